$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "Hoan thanh" -> "Hoàn thành"
Replace-Text "Hoan thanh Crawl." "Hoàn thành Crawl."

# 2. "traning" -> "training", "trung" -> "khung" (single run, no leading "- ")
Replace-Text "Thu thập tài liệu để làm tài liệu traning của ngành CS <trước hết là theo cái trung >" "Thu thập tài liệu để làm tài liệu training của ngành CS <trước hết là theo cái khung >"

# 3. "tóm" -> "Tóm" (capitalize)
Replace-Text "+ tóm tắt lại kiến thức có trong tài liệu. " "+ Tóm tắt lại kiến thức có trong tài liệu. "

# 4. "hoan thiện" -> "hoàn thiện"
Replace-Text " thế nào ,có nhiều cấp không ---> từ đó hoan thiện cái khung của mình." " thế nào ,có nhiều cấp không ---> từ đó hoàn thiện cái khung của mình."

# 5. "đẻ" -> "để", "dung" -> "dùng"
Replace-Text "- Xác định lại mục tiêu của đề tài ? phân loại để làm gì ? tại sao phải crawl  trong khi đó các thư viện có thể nó không cho download tài liệu về và các thư viện số đó đã phân loại rồi có cần thiết phải phân loại không ??????? <gợi ý là crawl các thông tin về bài báo như abtract đẻ phân loại làm dữ liệu làm giàu ontology dung để lưu vết các bài báo cung cấp các thông tin …… ,  >." "- Xác định lại mục tiêu của đề tài ? phân loại để làm gì ? tại sao phải crawl  trong khi đó các thư viện có thể nó không cho download tài liệu về và các thư viện số đó đã phân loại rồi có cần thiết phải phân loại không ??????? <gợi ý là crawl các thông tin về bài báo như abtract để phân loại làm dữ liệu làm giàu ontology dùng để lưu vết các bài báo cung cấp các thông tin …… ,  >."

# 6. "dưdợc" -> "được"
Replace-Text " Thông tin của 1 bài báo lấy dưdợc gì ? metadata,abtract, … " " Thông tin của 1 bài báo lấy được gì ? metadata,abtract, … "
